$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "hardware" is a duplicate of the english value in row 3 (B3); Excel's
# transliteration for it ("हार्डवेयरस") gets appended as a new row so the
# duplicate can be printed/reviewed.
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "hardware"
$ws.Range("C7").Value = "हार्डवेयरस"

# Match the Hindi-capable font used by the other transliteration/translation
# cells (C2:D6) so the new Devanagari text renders with the same style.
$ws.Range("C7").Font.Name = "Lohit Hindi"

# Leave the selection where the edit ended up.
$null = $ws.Range("D6").Select()
